# Generate Report for Handoff
# Adds a new handoff row (d60fcd1c-ad5c-4859-ae55-59a59c812f66) to every
# sheet of the localization-status workbook, mirroring the existing
# cf6ce968-bbd2-4664-ad2a-3218565f3e78 row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet - row 3
# ---------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

$ov.Range("A3").Value = "d60fcd1c-ad5c-4859-ae55-59a59c812f66.md"
$ov.Range("B3").Value = "Ready for handoff"
$ov.Range("C3").Value = "Ready for handoff"
$ov.Range("D3").Value = "2016-03-21 14:36:56"
$ov.Range("D3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$ov.Hyperlinks.Add($ov.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/46b442b4c7f8289d335235cf4dc2cc7420699bae/e2e/d60fcd1c-ad5c-4859-ae55-59a59c812f66.md", "", "", "d60fcd1c-ad5c-4859-ae55-59a59c812f66.md")
$ov.Range("A3").Font.Name = "Calibri"
$ov.Range("A3").Font.Size = 11
$ov.Range("A3").Font.Underline = $true
$ov.Range("A3").Font.Color = 15570276

# ---------------------------------------------------------------------
# zh-cn sheet - row 3
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("A3").Value = "d60fcd1c-ad5c-4859-ae55-59a59c812f66.md"
$zh.Range("B3").Value = ".md"
$zh.Range("C3").Value = "Ready for handoff"
$zh.Range("D3").Value = "d60fcd1c-ad5c-4859-ae55-59a59c812f66.efc061ea5a012367f66426a37431b07940d61a60.zh-cn.xlf"
$zh.Range("E3").Value = "2016-03-21 14:36:52"
$zh.Range("E3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$zh.Range("H3").Value = "0001-01-01 00:00:00"
$zh.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$zh.Range("J3").Value = "Include"

$zh.Hyperlinks.Add($zh.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/46b442b4c7f8289d335235cf4dc2cc7420699bae/e2e/d60fcd1c-ad5c-4859-ae55-59a59c812f66.md", "", "", "d60fcd1c-ad5c-4859-ae55-59a59c812f66.md")
$zh.Range("A3").Font.Name = "Calibri"
$zh.Range("A3").Font.Size = 11
$zh.Range("A3").Font.Underline = $true
$zh.Range("A3").Font.Color = 15570276

$zh.Hyperlinks.Add($zh.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5af2e2e6b124cac4fecb3d77294d6dbea7f10335/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/d60fcd1c-ad5c-4859-ae55-59a59c812f66.efc061ea5a012367f66426a37431b07940d61a60.zh-cn.xlf", "", "", "d60fcd1c-ad5c-4859-ae55-59a59c812f66.efc061ea5a012367f66426a37431b07940d61a60.zh-cn.xlf")
$zh.Range("D3").Font.Name = "Calibri"
$zh.Range("D3").Font.Size = 11
$zh.Range("D3").Font.Underline = $true
$zh.Range("D3").Font.Color = 15570276

# ---------------------------------------------------------------------
# de-de sheet - row 3
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("A3").Value = "d60fcd1c-ad5c-4859-ae55-59a59c812f66.md"
$de.Range("B3").Value = ".md"
$de.Range("C3").Value = "Ready for handoff"
$de.Range("D3").Value = "d60fcd1c-ad5c-4859-ae55-59a59c812f66.efc061ea5a012367f66426a37431b07940d61a60.de-de.xlf"
$de.Range("E3").Value = "2016-03-21 14:36:56"
$de.Range("E3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$de.Range("H3").Value = "0001-01-01 00:00:00"
$de.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$de.Range("J3").Value = "Include"

$de.Hyperlinks.Add($de.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/46b442b4c7f8289d335235cf4dc2cc7420699bae/e2e/d60fcd1c-ad5c-4859-ae55-59a59c812f66.md", "", "", "d60fcd1c-ad5c-4859-ae55-59a59c812f66.md")
$de.Range("A3").Font.Name = "Calibri"
$de.Range("A3").Font.Size = 11
$de.Range("A3").Font.Underline = $true
$de.Range("A3").Font.Color = 15570276

$de.Hyperlinks.Add($de.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/52fbf70eea2eebbc26d9d319f762103ee7823f1e/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/d60fcd1c-ad5c-4859-ae55-59a59c812f66.efc061ea5a012367f66426a37431b07940d61a60.de-de.xlf", "", "", "d60fcd1c-ad5c-4859-ae55-59a59c812f66.efc061ea5a012367f66426a37431b07940d61a60.de-de.xlf")
$de.Range("D3").Font.Name = "Calibri"
$de.Range("D3").Font.Size = 11
$de.Range("D3").Font.Underline = $true
$de.Range("D3").Font.Color = 15570276
